$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 103-105)
$rows = @(
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44832; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Especial"; M=400;  N=22000; O=23000; P=22500; Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=2250; T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44832; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Primera";  M=400;  N=19000; O=20000; P=19500; Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=1950; T=10 },
    @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44832; E=4; F="Fruta"; G=100107; H="Otros"; I=100107002; J="Chirimoya"; K="Cultivar IV Región"; L="Segunda";  M=360;  N=15000; O=16000; P=15500; Q="$/bandeja 10 kilos"; R="Provincia de Limarí"; S=1550; T=10 }
)

$startRow = 103
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value2 = $data.A
    $ws.Cells.Item($r, 2).Value2 = $data.B
    $ws.Cells.Item($r, 3).Value2 = $data.C

    $ws.Cells.Item($r, 4).Value2 = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value2 = $data.E
    $ws.Cells.Item($r, 6).Value2 = $data.F
    $ws.Cells.Item($r, 7).Value2 = $data.G
    $ws.Cells.Item($r, 8).Value2 = $data.H
    $ws.Cells.Item($r, 9).Value2 = $data.I
    $ws.Cells.Item($r, 10).Value2 = $data.J
    $ws.Cells.Item($r, 11).Value2 = $data.K
    $ws.Cells.Item($r, 12).Value2 = $data.L
    $ws.Cells.Item($r, 13).Value2 = $data.M
    $ws.Cells.Item($r, 14).Value2 = $data.N
    $ws.Cells.Item($r, 15).Value2 = $data.O
    $ws.Cells.Item($r, 16).Value2 = $data.P
    $ws.Cells.Item($r, 17).Value2 = $data.Q
    $ws.Cells.Item($r, 18).Value2 = $data.R
    $ws.Cells.Item($r, 19).Value2 = $data.S
    $ws.Cells.Item($r, 20).Value2 = $data.T
}

Write-Host "Done writing rows."
